$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("daily_report")
$ws.Activate()

# --- Row 24: fill in the day's actuals (E) and the updated projection (P) ---
$ws.Range("E24").Value = 9113

# The rest of row 24's formulas follow the same pattern used by the rows
# above them (the sheet is a day-over-day report where each new day's row
# gets the same formulas copied down) so recreate them explicitly.
$ws.Range("F24").Formula = "=E24-D24"
$ws.Range("G24").Formula = "=(E24-`$D`$2)/A24"
$ws.Range("H24").Formula = "=(E24/D24-1)*100"
$ws.Range("I24").Formula = "=(POWER((E24/`$D`$3),1/A24)-1)*100"
$ws.Range("K24").Formula = "=E24-J24"
$ws.Range("M24").Formula = "=L24+E24"
$ws.Range("N24").Formula = "=E24/`$D`$2*100"

$ws.Range("P24").Value = 87985.1

# --- Restore the last-used selection ---
$null = $ws.Range("P25").Select()
